$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the value to be stored as text even if it looks like a number
    # (matches the source data, which keeps these as inline/shared strings).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "61.696.88"
$ws.Range("E2").Value = "  +1.06%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.398.88"
$ws.Range("E3").Value = "  +0.20%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
Set-TextValue $ws.Range("D5") "577.37"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6
Set-TextValue $ws.Range("D6") "143.83"
$ws.Range("E6").Value = "  +0.93%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.474"
$ws.Range("E8").Value = "  -0.34%  "

# Row 9
Set-TextValue $ws.Range("D9") "7.61"
$ws.Range("E9").Value = "  -0.30%  "

# Row 10
$ws.Range("E10").Value = "  -0.61%  "

# Row 11
$ws.Range("E11").Value = "  -1.12%  "

# Row 12
Set-TextValue $ws.Range("D12") "3.978.46"
$ws.Range("E12").Value = "  +0.18%  "

# Row 13
$ws.Range("E13").Value = "  -0.19%  "

# Row 14
Set-TextValue $ws.Range("D14") "28.05"
$ws.Range("E14").Value = "  +0.91%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.385.84"
$ws.Range("E15").Value = "  -0.26%  "

# Row 16
$ws.Range("E16").Value = "  -0.56%  "

# Row 17
Set-TextValue $ws.Range("D17") "61.738.92"
$ws.Range("E17").Value = "  +0.95%  "

# Row 18
$ws.Range("E18").Value = "  +0.57%  "

# Row 19
Set-TextValue $ws.Range("D19") "13.69"
$ws.Range("E19").Value = "  +0.06%  "

# Row 20
$ws.Range("E20").Value = "  +1.87%  "

# Row 21
Set-TextValue $ws.Range("D21") "389.52"
$ws.Range("E21").Value = "  +1.55%  "

# Row 22
Set-TextValue $ws.Range("D22") "74.68"
$ws.Range("E22").Value = "  -0.38%  "

# Row 23
$ws.Range("E23").Value = "  -0.71%  "

# Row 24
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("E25").Value = "  -3.16%  "

# Row 26
Set-TextValue $ws.Range("D26") "0.181"
$ws.Range("E26").Value = "  +0.06%  "

# Row 27
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D27") "7.43"
$ws.Range("E27").Value = "  +1.06%  "

# Row 28
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D28") "0.997"
$ws.Range("E28").Value = "  -0.20%  "

# Row 29
$ws.Range("E29").Value = "  -0.32%  "

# Row 30
$ws.Range("E30").Value = "  -0.38%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.41"
$ws.Range("E31").Value = "  +0.05%  "

# Row 32
$ws.Range("E32").Value = "  -0.05%  "

# Row 33
Set-TextValue $ws.Range("D33") "23.42"
$ws.Range("E33").Value = "  +0.22%  "

# Row 34
$ws.Range("E34").Value = "  -0.51%  "

# Row 35
Set-TextValue $ws.Range("D35") "168.43"
$ws.Range("E35").Value = "  +1.09%  "

# Row 36
Set-TextValue $ws.Range("D36") "5.11"
$ws.Range("E36").Value = "  +1.34%  "

# Row 37
Set-TextValue $ws.Range("D37") "3.430.84"
$ws.Range("E37").Value = "  +0.23%  "

# Row 38
$ws.Range("E38").Value = "  +0.09%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.0764"
$ws.Range("E39").Value = "  -0.92%  "

# Row 40
Set-TextValue $ws.Range("D40") "27.09"
$ws.Range("E40").Value = "  -0.54%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.784"
$ws.Range("E41").Value = "  +0.44%  "

# Row 42
Set-TextValue $ws.Range("D42") "4.44"
$ws.Range("E42").Value = "  +0.92%  "

# Row 43
Set-TextValue $ws.Range("D43") "1.67"
$ws.Range("E43").Value = "  -0.11%  "

# Row 44
$ws.Range("E44").Value = "  +2.13%  "

# Row 45
Set-TextValue $ws.Range("D45") "2.480.39"
$ws.Range("E45").Value = "  +1.07%  "

# Row 46
Set-TextValue $ws.Range("D46") "22.82"
$ws.Range("E46").Value = "  -1.28%  "

# Row 47
$ws.Range("E47").Value = "  -1.07%  "

# Row 48
$ws.Range("E48").Value = "  +0.01%  "

# Row 49
$ws.Range("E49").Value = "  -0.38%  "

# Row 50
$ws.Range("E50").Value = "  -5.80%  "

# Row 51
$ws.Range("E51").Value = "  -1.13%  "
